# Updates cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51.
# Mirrors the commit "Updated cryptos list ... with GitHub Actions":
# plain text value replacements, no formatting/style changes intended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.618.37'
$ws.Range("E2").Value = '  +2.53%  '
$ws.Range("D3").Value = '2.629.75'
$ws.Range("E3").Value = '  +2.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("D9").Value = '2.629.00'
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("E12").Value = '  -5.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.368'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.77%  '
$ws.Range("D14").Value = '3.090.07'
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = '60.643.83'
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("E17").Value = '  +4.74%  '
$ws.Range("D18").Value = '2.627.20'
$ws.Range("E18").Value = '  +2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.28%  '
$ws.Range("E20").Value = '  +2.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.78%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("E24").Value = '  +9.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +1.84%  '
$ws.Range("E28").Value = '  +5.32%  '
$ws.Range("D29").Value = '0.0₃0794'
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("E30").Value = '  +10.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '162.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.08%  '
$ws.Range("E35").Value = '  +5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.980'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.74%  '
$ws.Range("E37").Value = '  +7.22%  '
$ws.Range("E38").Value = '  +8.89%  '
$ws.Range("E39").Value = '  +1.82%  '
$ws.Range("E40").Value = '  +6.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.847'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '300.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '134.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.97%  '
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.606'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.04%  '
$ws.Range("E49").Value = '  +3.48%  '
$ws.Range("E50").Value = '  +4.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.96'
$ws.Range("D51").Style = "Normal"
